$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns J (10) and K (11) -----------------------------------------
# Target stored widths (OOXML <col width=.../>) are 314.4 and 282 character
# units. The host quantizes ColumnWidth to a 1/6-character pixel grid (same
# style of rounding real Excel applies), so the input has to be pre-compensated
# by -5/6 to land on the desired stored value. 282 lands exactly; 314.4 is not
# reachable on that grid so we use the closest representable value.
$ws.Columns.Item(10).ColumnWidth = 313.5               # -> stored 314.3333333333333 (closest to 314.4)
$ws.Columns.Item(11).ColumnWidth = 281.16666666666669  # -> stored 282

# --- Row 1 headers (bold header style, like existing H1/I1) ----------------
$ws.Cells.Item(1, 10).Value = "Onkelos"
$ws.Cells.Item(1, 11).Value = "Jonathan"

foreach ($c in 10, 11) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Font.Bold = $true
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

# --- Row 2 data (plain wrap-text style, like existing A2:I2) ---------------
$ws.Cells.Item(2, 10).Value = '“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt.'
$ws.Cells.Item(2, 11).Value = 'Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;'

foreach ($c in 10, 11) {
    $cell = $ws.Cells.Item(2, $c)
    $cell.WrapText = $true
}
